# Applies the "Updated symbol list on Tue Feb 14 05:56:25 UTC 2023 with
# GitHub Actions" commit to cryptos.xlsx: the crypto-exchange-token rows
# (B6:E17) shifted down one slot to make room for GateToken at row 6, and
# the Price / Volume(1h) columns were refreshed with the latest scrape
# across most data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextCell {
    param($Sheet, $CellRef, $NewValue)
    # The source data keeps every cell (coin name, link, price, % change)
    # as literal text, even values that look numeric ("292.18") or like a
    # percentage ("-7.26%"). A leading apostrophe forces Excel to store the
    # value as text instead of inferring a number/date/percent, and resetting
    # the style back to "Normal" afterwards drops the quote-prefix formatting
    # Excel would otherwise tag the cell with, keeping the cell style untouched.
    $Sheet.Range($CellRef).Value = "'" + $NewValue
    $Sheet.Range($CellRef).Style = "Normal"
}

Set-TextCell $ws "D2" "292.18"
Set-TextCell $ws "E2" "-7.26%"
Set-TextCell $ws "D3" "40.34"
Set-TextCell $ws "E3" "-1.59%"
Set-TextCell $ws "E4" "-2.41%"
Set-TextCell $ws "D5" "0.07324"
Set-TextCell $ws "E5" "-3.68%"
Set-TextCell $ws "B6" "GateToken"
Set-TextCell $ws "C6" "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextCell $ws "D6" "4.299"
Set-TextCell $ws "E6" "-0.55%"
Set-TextCell $ws "B7" "FTXToken"
Set-TextCell $ws "C7" "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextCell $ws "D7" "1.535"
Set-TextCell $ws "E7" "-7.94%"
Set-TextCell $ws "B8" "MXToken"
Set-TextCell $ws "C8" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextCell $ws "D8" "0.9278"
Set-TextCell $ws "E8" "0.04%"
Set-TextCell $ws "B9" "BTSEToken"
Set-TextCell $ws "C9" "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextCell $ws "D9" "2.369"
Set-TextCell $ws "E9" "-2.27%"
Set-TextCell $ws "B10" "LiechtensteinCryptoassetsExchange"
Set-TextCell $ws "C10" "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextCell $ws "D10" "0.1177"
Set-TextCell $ws "E10" "-1.72%"
Set-TextCell $ws "B11" "WazirX"
Set-TextCell $ws "C11" "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextCell $ws "D11" "0.1737"
Set-TextCell $ws "E11" "-4.45%"
Set-TextCell $ws "B12" "BitrueCoin"
Set-TextCell $ws "C12" "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextCell $ws "D12" "0.04331"
Set-TextCell $ws "E12" "4.53%"
Set-TextCell $ws "B13" "MandalaExchangeToken"
Set-TextCell $ws "C13" "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextCell $ws "D13" "0.08630"
Set-TextCell $ws "E13" "-4.28%"
Set-TextCell $ws "B14" "BitMartToken"
Set-TextCell $ws "C14" "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextCell $ws "D14" "0.1056"
Set-TextCell $ws "E14" "0.19%"
Set-TextCell $ws "B15" "BitForexToken"
Set-TextCell $ws "C15" "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextCell $ws "D15" "0.001276"
Set-TextCell $ws "E15" "-0.89%"
Set-TextCell $ws "B16" "TigerCash"
Set-TextCell $ws "C16" "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextCell $ws "D16" "0.006030"
Set-TextCell $ws "E16" "3.87%"
Set-TextCell $ws "B17" "LEO"
Set-TextCell $ws "C17" "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell $ws "D17" "3.338"
Set-TextCell $ws "E17" "0.20%"
Set-TextCell $ws "E18" "-1.99%"
Set-TextCell $ws "D19" "7.974"
Set-TextCell $ws "E19" "5.12%"
Set-TextCell $ws "D20" "0.1401"
Set-TextCell $ws "E20" "3.62%"
Set-TextCell $ws "E21" "-2.21%"
Set-TextCell $ws "D22" "0.03946"
Set-TextCell $ws "E22" "-2.01%"
Set-TextCell $ws "E23" "-0.93%"
Set-TextCell $ws "D24" "0.003777"
Set-TextCell $ws "E24" "-7.08%"
Set-TextCell $ws "D25" "0.0001281"
Set-TextCell $ws "E25" "0.80%"
Set-TextCell $ws "D26" "0.0003728"
Set-TextCell $ws "E26" "-95.04%"
Set-TextCell $ws "D38" "0.02268"
Set-TextCell $ws "E38" "-6.24%"
Set-TextCell $ws "D39" "0.04981"
Set-TextCell $ws "E39" "-3.63%"
Set-TextCell $ws "D40" "0.005905"
Set-TextCell $ws "E40" "78.78%"
Set-TextCell $ws "D41" "0.007684"
Set-TextCell $ws "E41" "-0.71%"
Set-TextCell $ws "D42" "0.1284"
Set-TextCell $ws "E42" "-1.21%"
Set-TextCell $ws "D43" "0.007348"
Set-TextCell $ws "E43" "-3.37%"
Set-TextCell $ws "D44" "0.008284"
Set-TextCell $ws "E44" "-3.49%"
Set-TextCell $ws "D45" "0.2916"
Set-TextCell $ws "E45" "-14.58%"
Set-TextCell $ws "D46" "0.00006312"
Set-TextCell $ws "E46" "-4.23%"
Set-TextCell $ws "E47" "0.06%"
Set-TextCell $ws "E48" "-90.71%"
Set-TextCell $ws "E49" "0.06%"
Set-TextCell $ws "E50" "0.06%"

Write-Host "Applied 92 cell updates to Sheet1"
